$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "66.978.04"
$ws.Range("E2").Value = "  +4.38%  "

# Row 3
Set-TextValue $ws.Range("D3") "3.440.63"
$ws.Range("E3").Value = "  +3.67%  "

# Row 4
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
Set-TextValue $ws.Range("D5") "581.88"
$ws.Range("E5").Value = "  +5.43%  "

# Row 6
Set-TextValue $ws.Range("D6") "184.68"
$ws.Range("E6").Value = "  +6.61%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.632"
$ws.Range("E7").Value = "  +2.52%  "

# Row 8
Set-TextValue $ws.Range("D8") "3.434.09"
$ws.Range("E8").Value = "  +3.73%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.173"
$ws.Range("E10").Value = "  +1.26%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.647"
$ws.Range("E11").Value = "  +2.89%  "

# Row 12
Set-TextValue $ws.Range("D12") "56.32"
$ws.Range("E12").Value = "  +5.36%  "

# Row 13
$ws.Range("E13").Value = "  -1.02%  "

# Row 14
Set-TextValue $ws.Range("D14") "9.43"
$ws.Range("E14").Value = "  +4.52%  "

# Row 15
Set-TextValue $ws.Range("D15") "3.990.03"
$ws.Range("E15").Value = "  +3.72%  "

# Row 16
$ws.Range("E16").Value = "  +3.24%  "

# Row 17
Set-TextValue $ws.Range("D17") "3.433.00"
$ws.Range("E17").Value = "  +3.63%  "

# Row 18
$ws.Range("E18").Value = "  +0.30%  "

# Row 19
Set-TextValue $ws.Range("D19") "66.789.25"
$ws.Range("E19").Value = "  +2.10%  "

# Row 20
Set-TextValue $ws.Range("D20") "12.09"
$ws.Range("E20").Value = "  +3.61%  "

# Row 21
Set-TextValue $ws.Range("D21") "1.02"
$ws.Range("E21").Value = "  +3.39%  "

# Row 22
Set-TextValue $ws.Range("D22") "483.51"
$ws.Range("E22").Value = "  +7.07%  "

# Row 23
Set-TextValue $ws.Range("D23") "16.86"
$ws.Range("E23").Value = "  +22.23%  "

# Row 24
Set-TextValue $ws.Range("D24") "5.14"
$ws.Range("E24").Value = "  +3.17%  "

# Row 25
Set-TextValue $ws.Range("D25") "4.40"
$ws.Range("E25").Value = "  +8.14%  "

# Row 26
Set-TextValue $ws.Range("D26") "89.84"
$ws.Range("E26").Value = "  +3.75%  "

# Row 27
Set-TextValue $ws.Range("D27") "11.06"

# Row 28
Set-TextValue $ws.Range("D28") "2.95"
$ws.Range("E28").Value = "  +3.08%  "

# Row 29
Set-TextValue $ws.Range("D29") "9.16"
$ws.Range("E29").Value = "  +6.92%  "

# Row 30
Set-TextValue $ws.Range("D30") "31.33"
$ws.Range("E30").Value = "  +1.54%  "

# Row 31
Set-TextValue $ws.Range("D31") "7.16"
$ws.Range("E31").Value = "  +9.32%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D32") "11.72"
$ws.Range("E32").Value = "  +2.86%  "

# Row 33
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D33") "63.99"
$ws.Range("E33").Value = "  +5.53%  "

# Row 34
Set-TextValue $ws.Range("D34") "590.26"
$ws.Range("E34").Value = "  +4.12%  "

# Row 35
Set-TextValue $ws.Range("D35") "0.112"
$ws.Range("E35").Value = "  +5.25%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D36") "0.148"
$ws.Range("E36").Value = "  +5.80%  "

# Row 37
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D37") "1.00"
$ws.Range("E37").Value = "  +0.00%  "

# Row 38
$ws.Range("E38").Value = "  +0.80%  "

# Row 39
Set-TextValue $ws.Range("D39") "36.46"
$ws.Range("E39").Value = "  +3.93%  "

# Row 40
Set-TextValue $ws.Range("D40") "0.386"
$ws.Range("E40").Value = "  +5.47%  "

# Row 41
Set-TextValue $ws.Range("D41") "0.0₃0769"
$ws.Range("E41").Value = "  +4.93%  "

# Row 42
Set-TextValue $ws.Range("D42") "3.191.27"
$ws.Range("E42").Value = "  +4.61%  "

# Row 43
Set-TextValue $ws.Range("D43") "2.91"
$ws.Range("E43").Value = "  +6.20%  "

# Row 44
$ws.Range("E44").Value = "  +4.38%  "

# Row 45
$ws.Range("E45").Value = "  +4.96%  "

# Row 46
Set-TextValue $ws.Range("D46") "2.78"
$ws.Range("E46").Value = "  +22.49%  "

# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D47") "0.135"
$ws.Range("E47").Value = "  +1.58%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D48") "3.20"
$ws.Range("E48").Value = "  +1.47%  "

# Row 49
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D49") "8.72"
$ws.Range("E49").Value = "  +7.28%  "

# Row 50
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D50") "0.999"
$ws.Range("E50").Value = "  -0.06%  "

# Row 51
Set-TextValue $ws.Range("D51") "140.00"
$ws.Range("E51").Value = "  -1.26%  "
